# Generate Report for Handoff
# The handoff XLIFF for "6b3b1a50-7ade-428d-9c8f-0c2b58ec58a9.md" was
# (re)generated, so its "Latest Handoff Datetime" / "Latest HO Xliff
# Generate Date" timestamps are refreshed on all three sheets
# (Overview, zh-cn, de-de) for that file's row (row 4).

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G4").Value = "2016-10-27 07:50:44"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H4").Value = "2016-10-27 07:50:31"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H4").Value = "2016-10-27 07:50:44"
